$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 167, shifting existing rows 167:293 down to 168:294
$ws.Rows.Item(167).Insert()

# Populate the newly inserted row 167 with the new record
$ws.Cells.Item(167, 1).Value = 8
$ws.Cells.Item(167, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(167, 3).Value = "Coquimbo"
$ws.Cells.Item(167, 4).Value = 44977
$ws.Cells.Item(167, 4).NumberFormat = $ws.Cells.Item(168, 4).NumberFormat
$ws.Cells.Item(167, 5).Value = 4
$ws.Cells.Item(167, 6).Value = 100112037
$ws.Cells.Item(167, 7).Value = "Cebollín"
$ws.Cells.Item(167, 8).Value = "Sin especificar"
$ws.Cells.Item(167, 9).Value = "Primera"
$ws.Cells.Item(167, 10).Value = 1600
$ws.Cells.Item(167, 11).Value = 1000
$ws.Cells.Item(167, 12).Value = 1200
$ws.Cells.Item(167, 13).Value = 1100
$ws.Cells.Item(167, 14).Value = "`$/paquete 6 unidades"
$ws.Cells.Item(167, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(167, 16).Value = 183
$ws.Cells.Item(167, 17).Value = 6
$ws.Cells.Item(167, 18).Value = "Hortaliza"
